# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values for the df08d9e8-... file row on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-25 03:00:15"
$wsZhCn.Range("H3").Value = "2016-03-25 03:00:58"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-25 03:00:21"
$wsDeDe.Range("H3").Value = "2016-03-25 03:01:06"
